$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "HP Victus AMD Ryzen 7 Hexa Core 7445HS - (16 GB/512 GB SSD/Windows 11 Home/4 GB Graphics/NVIDIA GeForce RTX 2050) 15-fb3122AX Gaming Laptop (15.6 Inch, Performance Blue, Chrome Logo, 2.29 Kg, With MS Office)"
$ws.Range("B2").Value = 63990
$ws.Range("C2").Value = 4.4
$ws.Range("D2").Value = "AMD Ryzen 7 Hexa Core"
$ws.Range("E2").Value = "512 GB"
$ws.Range("F2").Value = [double]"6.87607438662291e-05"

$ws.Range("A3").Value = "Acer Aspire 7 Intel Core i5 13th Gen 13420H - (16 GB/512 GB SSD/Windows 11 Home/6 GB Graphics/NVIDIA GeForce RTX 3050) A715-79G Gaming Laptop (15.6 Inch, Black, 1.99 Kg)"
$ws.Range("B3").Value = 62990
$ws.Range("C3").Value = 4.3
$ws.Range("D3").Value = "Intel Core i5 13th Gen"
$ws.Range("E3").Value = "512 GB"
$ws.Range("F3").Value = [double]"6.826480393713287e-05"

$ws.Range("A4").Value = "Acer Aspire 7 Intel Core 5 210H - (16 GB/512 GB SSD/Windows 11 Home/6 GB Graphics/NVIDIA GeForce RTX 4050) A715-79G Gaming Laptop (15.6 Inch, Obsidian Black, 1.99 Kg)"
$ws.Range("B4").Value = 69990
$ws.Range("C4").Value = 4.3
$ws.Range("D4").Value = "Intel Core 5"
$ws.Range("E4").Value = "512 GB"
$ws.Range("F4").Value = [double]"6.143734819259894e-05"

$ws.Range("A5").Value = "ASUS TUF Gaming A15 (2025) AMD Ryzen 7 Hexa Core 7445HS - (16 GB/512 GB SSD/Windows 11 Home/4 GB Graphics/NVIDIA GeForce RTX 3050/144 Hz) FA506NCG-HN199W Gaming Laptop (15.6 Inch, Graphite Black, 2.3 Kg)"
$ws.Range("B5").Value = 72990
$ws.Range("C5").Value = 4.4
$ws.Range("D5").Value = "AMD Ryzen 7 Hexa Core"
$ws.Range("E5").Value = "512 GB"
$ws.Range("F5").Value = [double]"6.028223044252638e-05"

$ws.Range("A6").Value = "Lenovo LOQ Essential Intel Core i5 12th Gen 12450HX - (16 GB/512 GB SSD/Windows 11 Home/4 GB Graphics/NVIDIA GeForce RTX 3050A) 15IAX9E Gaming Laptop (15.6 Inch, Luna Grey, 1.77 kg, With MS Office)"
$ws.Range("B6").Value = 68980
$ws.Range("C6").Value = 3.9
$ws.Range("D6").Value = "Intel Core i5 12th Gen"
$ws.Range("E6").Value = "512 GB"
$ws.Range("F6").Value = [double]"5.65381269933314e-05"

$ws.Range("A7").Value = "Lenovo LOQ Intel Core i5 13th Gen 13450HX - (16 GB/512 GB SSD/Windows 11 Home/6 GB Graphics/NVIDIA GeForce RTX 3050) LOQ 15IRX9D2 Gaming Laptop (15.6 Inch, Luna Grey, 2.38 kg, With MS Office)"
$ws.Range("B7").Value = 78990
$ws.Range("C7").Value = 4.3
$ws.Range("D7").Value = "Intel Core i5 13th Gen"
$ws.Range("E7").Value = "512 GB"
$ws.Range("F7").Value = [double]"5.443727054057475e-05"

$ws.Range("A8").Value = "DELL G15 Intel Core i5 13th Gen 13450HX - (16 GB/512 GB SSD/Windows 11 Home/6 GB Graphics/NVIDIA GeForce RTX 3050/120 Hz) 5530 Gaming Laptop (15.6 Inch, Dark Shadow Gray With Black Thermal Shelf, 2.65 Kg, With MS Office)"
$ws.Range("B8").Value = 79990
$ws.Range("C8").Value = 4.2
$ws.Range("D8").Value = "Intel Core i5 13th Gen"
$ws.Range("E8").Value = "512 GB"
$ws.Range("F8").Value = [double]"5.250656332041506e-05"

$ws.Range("A9").Value = "Acer NITRO V 16S Intel Core 5 - (16 GB/512 GB SSD/Windows 11 Home/8 GB Graphics/NVIDIA GeForce RTX NVIDIA GeForce RTX 5050/180 Hz) ANV16S-71 Gaming Laptop (16 Inch, Obsidian Black, 2.1 kg)"
$ws.Range("B9").Value = 94990
$ws.Range("C9").Value = 4.8
$ws.Range("D9").Value = "Intel Core 5"
$ws.Range("E9").Value = "512 GB"
$ws.Range("F9").Value = [double]"5.053163490893778e-05"

$ws.Range("A10").Value = "Lenovo LOQ Intel Core i5 13th Gen 13450HX - (16 GB/512 GB SSD/Windows 11 Home/6 GB Graphics/NVIDIA GeForce RTX 4050) 15IRX9 Gaming Laptop (15.6 Inch, Luna Grey, 2.38 Kg, With MS Office)"
$ws.Range("B10").Value = 87990
$ws.Range("C10").Value = 4.4
$ws.Range("D10").Value = "Intel Core i5 13th Gen"
$ws.Range("E10").Value = "512 GB"
$ws.Range("F10").Value = [double]"5.000568246391636e-05"

$ws.Range("A11").Value = "Acer Aspire 7 Intel Core i7 13th Gen 13620H - (16 GB/512 GB SSD/Windows 11 Home/6 GB Graphics/NVIDIA GeForce RTX 4050) A715-79G Gaming Laptop (15.6 Inch, Obsidian Black, 1.99 Kg)"
$ws.Range("B11").Value = 84990
$ws.Range("C11").Value = 4.2
$ws.Range("D11").Value = "Intel Core i7 13th Gen"
$ws.Range("E11").Value = "512 GB"
$ws.Range("F11").Value = [double]"4.941757853865161e-05"

$ws.Range("A12").Value = "HP Victus AMD Ryzen 7 Octa Core 260 - (24 GB/1 TB SSD/Windows 11 Home/8 GB Graphics/NVIDIA GeForce RTX 5050) 15-fb3185AX Gaming Laptop (15.6 Inch, Mica Silver, Black Chrome Logo, 2.29 Kg, With MS Office)"
$ws.Range("B12").Value = 99990
$ws.Range("C12").Value = 4.5
$ws.Range("D12").Value = "AMD Ryzen 7 Octa Core"
$ws.Range("E12").Value = "1 TB"
$ws.Range("F12").Value = [double]"4.5004500450045e-05"

$ws.Range("A13").Value = "Lenovo LOQ 2025 Intel Core i7 13th Gen 13700HX - (16 GB/1 TB SSD/Windows 11 Home/8 GB Graphics/NVIDIA GeForce RTX 5060) LoQ 15IRX10 Gaming Laptop (15.6 Inch, Luna Grey, 2.4 Kg, With MS Office)"
$ws.Range("B13").Value = 127990
$ws.Range("C13").Value = 4.5
$ws.Range("D13").Value = "Intel Core i7 13th Gen"
$ws.Range("E13").Value = "1 TB"
$ws.Range("F13").Value = [double]"3.515899679662474e-05"

$ws.Range("A14:F18").Delete() | Out-Null
